$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Asalto Magnético" card's description (D5) with the new
# "Coloca esta trampa..." trap text, and mark the cell as wrapped +
# vertically centered (matches the new cellXfs entry).
$newText = "Coloca esta trampa en un casillero desocupado. Bloquea de sus 8 direcciones una por cada nivel que tengas. Al ser pisado teletransporta al objetivo, trampa o proyectil a un casillero al azar: d8 dirección d6 distancia (cuenta en diagonal).`nPuedes seleccionar esta trampa para colocar otras trampas, efectos o habilidades no melé y estas se teletransportarán también.`nAl final de tu turno puedes modificar los bloqueos. Xendra."

$d5 = $ws.Range("D5")
$d5.Value = $newText
$d5.WrapText = $true
$d5.VerticalAlignment = -4108

# The card's description got a lot longer, so the row needs to grow to fit it.
$ws.Rows.Item(5).RowHeight = 375

# Selection moved from G5 to D5 (and the view scrolled so row 5 is visible).
$d5.Select()
